$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Interactive Presentation"
$ws.Range("B2").Value = "Interactive Presentation is an ideal instructional method for this course as it allows learners to engage actively with the content related to generative AI and storytelling. By incorporating interactive elements, such as real-time AI tool demonstrations and audience participation, learners can better understand complex concepts like AI-generated script ideation and ethical considerations. This method fosters an engaging learning environment where learners can ask questions and receive immediate feedback, enhancing their grasp of narrative structures and AI applications in media."

$ws.Range("A3").Value = "Brainstorming"
$ws.Range("B3").Value = "Brainstorming is particularly suitable for this course as it encourages creative thinking and collaboration among learners. When developing compelling script elements or identifying effective prompt terms for AI-generated storyboards, brainstorming sessions allow learners to explore diverse ideas and perspectives. This method supports the iterative process of refining AI-generated content by enabling learners to collectively generate innovative solutions and enhance their understanding of narrative components and visual storytelling techniques."

$ws.Range("A4").Value = "Demonstration"
$ws.Range("B4").Value = "Demonstration is a crucial instructional method for this course, as it provides learners with a visual and hands-on approach to understanding generative AI tools and techniques. By demonstrating how to use AI tools for script development, storyboarding, and video creation, learners can see the practical application of theoretical concepts. This method helps learners grasp complex processes, such as refining video scripts for clarity and tone, by breaking them down into manageable steps and showcasing real-time results, thereby enhancing their learning experience."

$ws.Range("A5").Value = "Drill and Practice"
$ws.Range("B5").Value = "Drill and Practice is an effective instructional method for reinforcing learners' understanding of AI-generated content and ethical considerations. By providing opportunities for learners to practice applying AI tools and techniques in a controlled environment, they can refine their skills and knowledge. This method helps learners build confidence in their ability to analyze AI outputs for bias and copyright risks, apply corrective actions, and develop best practices for ethical AI content generation. Through repeated practice, learners can better prepare for real-world challenges in media production."

$wb.Save()
